$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-9 down to 4-10
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the "M_full_png_remote_to_local_from_PI" measurement
$ws.Range("A3").Value = "M_full_png_remote_to_local_from_PI"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = 2.8386999999999998
$ws.Range("C3").Value = 2.7683
$ws.Range("D3").Formula = "= B3 - C3"

# Update column widths (A widened, B narrowed slightly)
$ws.Columns.Item(1).ColumnWidth = 32
$ws.Columns.Item(2).ColumnWidth = 18

# Update the active selection
$ws.Range("B5").Select()
